$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Grow the Table1 listobject from A8:K163 to A8:K175 (12 new rows) so the
#    new leave-card entries below belong to the table.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A8:K175"))

# ---------------------------------------------------------------------------
# 2) Row 161 ("FL(5-0-0)" -> "SL(104)") - clear the old 5-day undertime entry
#    and record the new SL particulars / days / remarks.
# ---------------------------------------------------------------------------
$ws.Range("B161").Value = "SL(104)"
$ws.Range("D161").ClearContents()
$ws.Range("H161").Value = 104
$ws.Range("K161").Value = "07/26/2023-12/29/2023"

# ---------------------------------------------------------------------------
# 3) Row 163 becomes the "2024" year-divider row (matches the look of the
#    existing "2023"/"2013" divider rows - bold, quote-prefixed text). Done
#    before row 162's own edits below so new shared-string entries line up
#    in the same order the original author created them in.
# ---------------------------------------------------------------------------
$ws.Range("A149").Copy()
$ws.Range("A163").PasteSpecial(-4122)
$ws.Range("A163").Value = "'2024"
$ws.Range("G163").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Row 162 ("SL(110-0-0)" -> "SL(70)") - new days value, updated H value and
#    remarks, plus the "EARNED " helper formula that every data row carries.
# ---------------------------------------------------------------------------
$ws.Range("B162").Value = "SL(70)"
$ws.Range("D162").Value = 12
$ws.Range("G162").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("H162").Value = 58
$ws.Range("K162").Value = "1/2023 - 7/11/2023"

# ---------------------------------------------------------------------------
# 5) Row 164 - first entry of 2024 ("SL(22-0-0)"), formatted like the normal
#    data rows (template: row 163's original blank styling).
# ---------------------------------------------------------------------------
$ws.Range("A163:K163").Copy()
$ws.Range("A164:K164").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A164").Formula = "=EDATE(A161,1)"
$ws.Range("B164").Value = "SL(22-0-0)"
$ws.Range("D164").Value = 8
$ws.Range("J164").Value = 14
$ws.Range("K164").Value = "1/2-31/2024"

# ---------------------------------------------------------------------------
# 6) Row 165 - blank continuation row, but B/D/H keep the "post-divider"
#    styling used right after a year-divider row (copied from row 149).
# ---------------------------------------------------------------------------
$ws.Range("A163:K163").Copy()
$ws.Range("A165:K165").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B149").Copy()
$ws.Range("B165").PasteSpecial(-4122)
$ws.Range("D149").Copy()
$ws.Range("D165").PasteSpecial(-4122)
$ws.Range("H149").Copy()
$ws.Range("H165").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A165").Formula = "=EDATE(A164,1)"

# ---------------------------------------------------------------------------
# 7) Rows 166-175 - ten blank placeholder rows (one per month), each one
#    carrying the standard data-row styling and the "EARNED " helper formula.
# ---------------------------------------------------------------------------
$ws.Range("A148:K148").Copy()
$ws.Range("A166:K175").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$prev = "A165"
for ($r = 166; $r -le 175; $r++) {
    $ws.Range("A$r").Formula = "=EDATE($prev,1)"
    $ws.Range("G$r").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
    $prev = "A$r"
}

# ---------------------------------------------------------------------------
# 8) Leave the cursor where the user left it on save.
# ---------------------------------------------------------------------------
$ws.Range("J163").Select()

$wb.Application.Calculate()
